# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    3  = 3125
    7  = 1666
    8  = 1621
    18 = 4
    21 = 46
    22 = 19
    23 = 370
    25 = 98
    28 = 88
    29 = 2100
    34 = 570
    37 = 341
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
